# Milestone 3 Presentation fixes
#
# 1) Slide 7 (shape "Google Shape;182;p19"): the 3rd paragraph's run text was
#    cut short - append the missing tail so it reads as a complete sentence.
# 2) Theme colors: the deck's primary (slide-facing) theme palette should use
#    the original "Default" color set instead of the "Macmorris" colors that
#    had been swapped in.

$p = $ppt.ActivePresentation

# --- 1) Fix the truncated sentence on slide 7 -----------------------------
$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(3)
$textRange = $shape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(3, 1)
$run = $paragraph.Runs(1, 1)
$run.Text = "Yesterday as a group we all sat down, merged our code, and ensured the code was clean and that the code followed some of the formalities we went over in class."

# --- 2) Restore the "Default" color scheme on the presentation's theme ---
$colorScheme = $p.Slides.Item(1).ThemeColorScheme
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 5800213    # dk2      158158
$colorScheme.Item(4).RGB  = 15987699   # lt2      F3F3F3
$colorScheme.Item(5).RGB  = 13077765   # accent1  058DC7
$colorScheme.Item(6).RGB  = 3322960    # accent2  50B432
$colorScheme.Item(7).RGB  = 1791725    # accent3  ED561B
$colorScheme.Item(8).RGB  = 61421      # accent4  EDEF00
$colorScheme.Item(9).RGB  = 15059748   # accent5  24CBE5
$colorScheme.Item(10).RGB = 7529828    # accent6  64E572
$colorScheme.Item(11).RGB = 13369378   # hlink    2200CC
$colorScheme.Item(12).RGB = 9116245    # folHlink 551A8B
